# Apply the updated cryptocurrency price/volume figures.
# Column D (Price) values that look numeric get auto-converted to
# floating point by Excel's normal type inference, which would lose
# the original plain-text representation (and introduce float noise).
# Force text entry via a temporary "@" (Text) number format, then
# clear the format again afterwards so the cell ends up with no
# explicit style, just like the source file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.196.09'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.660.55'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.26'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5214'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06321'
$ws.Range('D9').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.12'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07729'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.659.45'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.430'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.886.55'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5457'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.34%  '
$ws.Range('E16').Value = '  -2.70%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.94'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '26.244.84'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.663'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '193.18'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('E22').Value = '  -2.56%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.098'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.61%  '
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '138.51'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1243'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.34%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.213'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '16.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.06008'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.95%  '
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.575'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.330'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.12%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.649'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.21%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9817'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.779'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.409'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5906'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01586'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.944'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8632'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.003'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.040.29'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '99.63'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.87%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.801.16'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.95%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0₈109'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '57.09'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.007'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.120'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.44%  '
